# dispatch suite & order creation optimization
#
# - Normalize SKU codes on the "SKU_Data" sheet to uppercase ("2ss" -> "2SS",
#   "7ss" -> "7SS")
# - Tighten row height on the "Credentials" sheet's trailing blank row
# - Switch the active/selected sheet & selection from "Credentials" to
#   "SKU_Data" (with cell B5 selected there instead of G6)

$wb = $excel.ActiveWorkbook

$credentials = $wb.Worksheets.Item("Credentials")
$skuData = $wb.Worksheets.Item("SKU_Data")

# Uppercase the SKU code strings used on the SKU_Data sheet
$skuData.Range("A2").Value = "2SS"
$skuData.Range("A3").Value = "7SS"

# Shrink the empty row 5 on Credentials from 13 to 12.8 points
$credentials.Rows.Item(5).RowHeight = 12.8

# Make SKU_Data the active sheet with B5 selected (previously Credentials was
# active and SKU_Data had G6 selected)
$skuData.Range("B5").Select()
$skuData.Activate()
